$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "31.071.22"
$ws.Range("E2").Value = "  +1.20%  "
$ws.Range("D3").Value = "1.954.57"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.09"
$ws.Range("E5").Value = "  -0.42%  "
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("E7").Value = "  +1.34%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2974"
$ws.Range("E8").Value = "  +1.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06842"
$ws.Range("E9").Value = "  +0.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.16"
$ws.Range("E10").Value = "  -1.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "107.70"
$ws.Range("E11").Value = "  -4.29%  "
$ws.Range("D12").Value = "1.941.95"
$ws.Range("E13").Value = "  +1.08%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.457"
$ws.Range("E14").Value = "  -1.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.7080"
$ws.Range("E15").Value = "  +2.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "282.42"
$ws.Range("E16").Value = "  -4.63%  "
$ws.Range("D17").Value = "30.897.62"
$ws.Range("E17").Value = "  +0.43%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.29"
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007753"
$ws.Range("E19").Value = "  +0.52%  "
$ws.Range("B20").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D20").Value = "2.196.58"
$ws.Range("E20").Value = "  -0.22%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9999"
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.498"
$ws.Range("E22").Value = "  -3.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.003"
$ws.Range("E23").Value = "  +0.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.498"
$ws.Range("E24").Value = "  -1.63%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.826"
$ws.Range("E25").Value = "  -0.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.51"
$ws.Range("E26").Value = "  +0.70%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.04"
$ws.Range("E27").Value = "  -1.46%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.209"
$ws.Range("E28").Value = "  +0.72%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1057"
$ws.Range("E29").Value = "  -2.99%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.424"
$ws.Range("E30").Value = "  -1.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.583"
$ws.Range("E31").Value = "  -0.56%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.568"
$ws.Range("E32").Value = "  -3.64%  "
$ws.Range("E33").Value = "  -0.89%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04966"
$ws.Range("E34").Value = "  -2.51%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7575"
$ws.Range("E35").Value = "  -2.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.183"
$ws.Range("E36").Value = "  +1.86%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.729"
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02037"
$ws.Range("E38").Value = "  -2.53%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.706"
$ws.Range("E39").Value = "  +0.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.177"
$ws.Range("E40").Value = "  +5.95%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.493"
$ws.Range("E41").Value = "  +9.67%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "74.57"
$ws.Range("E42").Value = "  +6.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4506"
$ws.Range("E43").Value = "  +0.83%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "109.31"
$ws.Range("E44").Value = "  -1.94%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8820"
$ws.Range("E45").Value = "  +0.91%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.190"
$ws.Range("E46").Value = "  +11.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.000"
$ws.Range("E47").Value = "  -0.27%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "977.98"
$ws.Range("E48").Value = "  +8.01%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.444"
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("E50").Value = "  +0.92%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.2580"
$ws.Range("E51").Value = "  +1.83%  "
